# Weekly data update for "Fruta, Feria Lagunitas de Puerto Montt - Naranja"
# Insert a new week's record as row 749, shifting all following rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(749).Insert()

$ws.Cells.Item(749, 1).Value = 4
$ws.Cells.Item(749, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(749, 3).Value = "Los Lagos"
$ws.Cells.Item(749, 4).Value = 45075
$ws.Cells.Item(749, 5).Value = 10
$ws.Cells.Item(749, 6).Value = "Fruta"
$ws.Cells.Item(749, 7).Value = 100102
$ws.Cells.Item(749, 8).Value = "Cítricos"
$ws.Cells.Item(749, 9).Value = 100102005
$ws.Cells.Item(749, 10).Value = "Naranja"
$ws.Cells.Item(749, 11).Value = "Fukumoto"
$ws.Cells.Item(749, 12).Value = "Primera"
$ws.Cells.Item(749, 13).Value = 400
$ws.Cells.Item(749, 14).Value = 19000
$ws.Cells.Item(749, 15).Value = 20000
$ws.Cells.Item(749, 16).Value = 19500
$ws.Cells.Item(749, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(749, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(749, 19).Value = 1300
$ws.Cells.Item(749, 20).Value = 15
